$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows at position 14, shifting existing rows 14+ down to 16+
$ws.Rows("14:15").Insert()

# Apply the same style (s="2") used by neighboring data rows to the new cells
$ws.Range("B14:E15").Style = $ws.Range("B16:E16").Style

# Populate the new "expand" row (row 14)
$ws.Range("B14").Value = "expand"
$ws.Range("C14").Value = "Expand"
$ws.Range("D14").Value = "Uitklappen"
$ws.Range("E14").Value = "Maximieren"

# Populate the new "collapse" row (row 15)
$ws.Range("B15").Value = "collapse"
$ws.Range("C15").Value = "Collapse"
$ws.Range("D15").Value = "Inklappen"
$ws.Range("E15").Value = "Minimieren"

# Update view state: clear frozen/scrolled topLeftCell, set new selection
$ws.Activate()
$ws.Range("C9").Select()
